# =====================================================================
# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund holdings detail) directly
# before the existing "2022-Q3" worksheet, and adds the corresponding
# summary row on the "总计" (totals) sheet.
# =====================================================================

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ3    = $wb.Worksheets.Item(2)   # "2022-Q3" (existing fund list)

# ---------------------------------------------------------------
# 1) New "2022-Q4" sheet, inserted right before "2022-Q3"
# ---------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# Header row, styled like the "总计" header (bold + border)
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"
$wsTotal.Range("B1").Copy() | Out-Null
$wsQ4.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data rows (index col A + fund code/name/size/position/weight/value/rank)
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").Value = "'160425"
$wsQ4.Range("C2").Value = "华安创业板两年定期开放混合"
$wsQ4.Range("D2").Value = "'1.71"
$wsQ4.Range("E2").Value = "'97.00"
$wsQ4.Range("F2").Value = "'9.63"
$wsQ4.Range("G2").Value = "'0.1647"
$wsQ4.Range("H2").Value = 1
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").Value = "'005457"
$wsQ4.Range("C3").Value = "景顺长城量化小盘股票"
$wsQ4.Range("D3").Value = "'5.08"
$wsQ4.Range("E3").Value = "'94.36"
$wsQ4.Range("F3").Value = "'1.35"
$wsQ4.Range("G3").Value = "'0.0686"
$wsQ4.Range("H3").Value = 8
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("B4").Value = "'007047"
$wsQ4.Range("C4").Value = "长城核心优势混合"
$wsQ4.Range("D4").Value = "'1.37"
$wsQ4.Range("E4").Value = "'76.80"
$wsQ4.Range("F4").Value = "'3.32"
$wsQ4.Range("G4").Value = "'0.0455"
$wsQ4.Range("H4").Value = 2
$wsQ4.Range("A5").Value = 3
$wsQ4.Range("B5").Value = "'010797"
$wsQ4.Range("C5").Value = "长城优选回报六个月持有期混合A"
$wsQ4.Range("D5").Value = "'2.72"
$wsQ4.Range("E5").Value = "'31.00"
$wsQ4.Range("F5").Value = "'1.18"
$wsQ4.Range("G5").Value = "'0.0321"
$wsQ4.Range("H5").Value = 6
$wsQ4.Range("A6").Value = 4
$wsQ4.Range("B6").Value = "'200001"
$wsQ4.Range("C6").Value = "长城久恒灵活配置混合"
$wsQ4.Range("D6").Value = "'0.85"
$wsQ4.Range("E6").Value = "'94.35"
$wsQ4.Range("F6").Value = "'3.00"
$wsQ4.Range("G6").Value = "'0.0255"
$wsQ4.Range("H6").Value = 4
$wsQ4.Range("A7").Value = 5
$wsQ4.Range("B7").Value = "'008851"
$wsQ4.Range("C7").Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$wsQ4.Range("D7").Value = "'2.37"
$wsQ4.Range("E7").Value = "'71.22"
$wsQ4.Range("F7").Value = "'1.06"
$wsQ4.Range("G7").Value = "'0.0251"
$wsQ4.Range("H7").Value = 7
$wsQ4.Range("A8").Value = 6
$wsQ4.Range("B8").Value = "'002703"
$wsQ4.Range("C8").Value = "长城久源灵活配置混合A"
$wsQ4.Range("D8").Value = "'0.70"
$wsQ4.Range("E8").Value = "'94.29"
$wsQ4.Range("F8").Value = "'3.49"
$wsQ4.Range("G8").Value = "'0.0244"
$wsQ4.Range("H8").Value = 10
$wsQ4.Range("A9").Value = 7
$wsQ4.Range("B9").Value = "'015496"
$wsQ4.Range("C9").Value = "景顺中证1000指数增强C"
$wsQ4.Range("D9").Value = "'0.86"
$wsQ4.Range("E9").Value = "'92.30"
$wsQ4.Range("F9").Value = "'1.52"
$wsQ4.Range("G9").Value = "'0.0131"
$wsQ4.Range("H9").Value = 4
$wsQ4.Range("A10").Value = 8
$wsQ4.Range("B10").Value = "'015495"
$wsQ4.Range("C10").Value = "景顺中证1000指数增强A"
$wsQ4.Range("D10").Value = "'0.67"
$wsQ4.Range("E10").Value = "'92.30"
$wsQ4.Range("F10").Value = "'1.52"
$wsQ4.Range("G10").Value = "'0.0102"
$wsQ4.Range("H10").Value = 4
$wsQ4.Range("A11").Value = 9
$wsQ4.Range("B11").Value = "'010798"
$wsQ4.Range("C11").Value = "长城优选回报六个月持有期混合C"
$wsQ4.Range("D11").Value = "'0.41"
$wsQ4.Range("E11").Value = "'31.00"
$wsQ4.Range("F11").Value = "'1.18"
$wsQ4.Range("G11").Value = "'0.0048"
$wsQ4.Range("H11").Value = 6
$wsQ4.Range("A12").Value = 10
$wsQ4.Range("B12").Value = "'014381"
$wsQ4.Range("C12").Value = "长城久源灵活配置混合C"
$wsQ4.Range("D12").Value = "'0.06"
$wsQ4.Range("E12").Value = "'94.29"
$wsQ4.Range("F12").Value = "'3.49"
$wsQ4.Range("G12").Value = "'0.0021"
$wsQ4.Range("H12").Value = 10

# Match the "总计" A-column style (bold + border) down column A
$wsTotal.Range("A2").Copy() | Out-Null
$wsQ4.Range("A2:A12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q4 summary row above 2022-Q3
# ---------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# New row 2 inherits row-above formatting from Insert(); start clean
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 11
$wsTotal.Range("D2").Value = 0.42

# Re-apply the index-column style to A2 (to match A3 below)
$wsTotal.Range("A3").Copy() | Out-Null
$wsTotal.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The old 2022-Q3 row shifted down to row 3; its index becomes 1
$wsTotal.Range("A3").Value = 1

